# Apply scheduled market-data refresh updates to profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 512.625
$ws.Range("I33").Value = 458.66666
$ws.Range("K33").Value = 458.66666
$ws.Range("M33").Value = -229.66666
$ws.Range("H43").Value = 5506.5
$ws.Range("J43").Value = 4350.1665
$ws.Range("L43").Value = 4350.1665
$ws.Range("N43").Value = -4488.1665
$ws.Range("H63").Value = 99999
$ws.Range("J63").Value = 99999
$ws.Range("L63").Value = 99999
$ws.Range("N63").Value = -101247
$ws.Range("H66").Value = 99999
$ws.Range("J66").Value = 99999
$ws.Range("L66").Value = 299997
$ws.Range("N66").Value = -306237
$ws.Range("H70").Value = 763257.4399999999
$ws.Range("I70").Value = 2033309.9
$ws.Range("J70").Value = 1226
$ws.Range("K70").Value = 6099929.699999999
$ws.Range("L70").Value = 3678
$ws.Range("M70").Value = -6099659.699999999
$ws.Range("N70").Value = -4218
$ws.Range("H73").Value = 763257.4399999999
$ws.Range("I73").Value = 2033309.9
$ws.Range("J73").Value = 1226
$ws.Range("K73").Value = 6099929.699999999
$ws.Range("L73").Value = 3678
$ws.Range("M73").Value = -6098993.699999999
$ws.Range("N73").Value = -5550

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 279.76923
$ws.Range("I5").Value = 221.54546
$ws.Range("J5").Value = 600
$ws.Range("K5").Value = 221.54546
$ws.Range("L5").Value = 600
$ws.Range("M5").Value = -109.54546
$ws.Range("N5").Value = -824
$ws.Range("H97").Value = 2125.6667
$ws.Range("I97").Value = 799.6667
$ws.Range("J97").Value = 3451.6667
$ws.Range("K97").Value = 799.6667
$ws.Range("L97").Value = 3451.6667
$ws.Range("M97").Value = -303.6667
$ws.Range("N97").Value = -4443.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 279.76923
$ws.Range("I4").Value = 221.54546
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 221.54546
$ws.Range("L4").Value = 600
$ws.Range("M4").Value = -106.54546
$ws.Range("N4").Value = -830
$ws.Range("H22").Value = 976.125
$ws.Range("J22").Value = 1997.5
$ws.Range("L22").Value = 1997.5
$ws.Range("N22").Value = -2343.5
$ws.Range("H64").Value = 994.4
$ws.Range("J64").Value = 994.4
$ws.Range("L64").Value = 994.4
$ws.Range("N64").Value = -1444.4
$ws.Range("H67").Value = 994.4
$ws.Range("J67").Value = 994.4
$ws.Range("L67").Value = 994.4
$ws.Range("N67").Value = -2554.4
$ws.Range("H110").Value = 119497.5
$ws.Range("J110").Value = 119497.5
$ws.Range("L110").Value = 119497.5
$ws.Range("N110").Value = -127677.5
$ws.Range("H111").Value = 20000
$ws.Range("J111").Value = 20000
$ws.Range("L111").Value = 20000
$ws.Range("N111").Value = -28180
$ws.Range("H129").Value = 100000
$ws.Range("J129").Value = 100000
$ws.Range("L129").Value = 100000
$ws.Range("N129").Value = -110000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 501
$ws.Range("I22").Value = 501
$ws.Range("K22").Value = 501
$ws.Range("M22").Value = -151
$ws.Range("H58").Value = 2997.48
$ws.Range("I58").Value = 2501.2632
$ws.Range("K58").Value = 2501.2632
$ws.Range("M58").Value = -2298.2632
$ws.Range("H122").Value = 2862.5833
$ws.Range("I122").Value = 3055
$ws.Range("K122").Value = 9165
$ws.Range("M122").Value = -6715
$ws.Range("H136").Value = 2997.48
$ws.Range("I136").Value = 2501.2632
$ws.Range("K136").Value = 7503.7896
$ws.Range("M136").Value = -4953.7896

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 750
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 750
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 2250
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -3314

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 98534
$ws.Range("J103").Value = 98534
$ws.Range("L103").Value = 98534
$ws.Range("N103").Value = -100878
$ws.Range("H141").Value = 132994.67
$ws.Range("I141").Value = 105500
$ws.Range("J141").Value = 187984
$ws.Range("K141").Value = 105500
$ws.Range("L141").Value = 187984
$ws.Range("M141").Value = -100320
$ws.Range("N141").Value = -198344

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9151.08
$ws.Range("I7").Value = 9132.388999999999
$ws.Range("J7").Value = 9199.143
$ws.Range("K7").Value = 9132.388999999999
$ws.Range("L7").Value = 9199.143
$ws.Range("M7").Value = -9020.388999999999
$ws.Range("N7").Value = -9423.143
$ws.Range("H22").Value = 3333
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 3333
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 3333
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -3923
$ws.Range("H27").Value = 3333
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 3333
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 3333
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -3547
$ws.Range("H46").Value = 1780.2858
$ws.Range("I46").Value = 1508.5
$ws.Range("J46").Value = 1889
$ws.Range("K46").Value = 1508.5
$ws.Range("L46").Value = 1889
$ws.Range("M46").Value = -1320.5
$ws.Range("N46").Value = -2265
$ws.Range("H63").Value = 66329.664
$ws.Range("J63").Value = 66329.664
$ws.Range("L63").Value = 66329.664
$ws.Range("N63").Value = -67827.664
$ws.Range("H66").Value = 66329.664
$ws.Range("J66").Value = 66329.664
$ws.Range("L66").Value = 198988.992
$ws.Range("N66").Value = -206476.992
$ws.Range("H100").Value = 19254156
$ws.Range("I100").Value = 3492
$ws.Range("J100").Value = 50055220
$ws.Range("K100").Value = 3492
$ws.Range("L100").Value = 50055220
$ws.Range("M100").Value = -2951
$ws.Range("N100").Value = -50056302
$ws.Range("H126").Value = 9151.08
$ws.Range("I126").Value = 9132.388999999999
$ws.Range("J126").Value = 9199.143
$ws.Range("K126").Value = 27397.167
$ws.Range("L126").Value = 27597.429
$ws.Range("M126").Value = -24927.167
$ws.Range("N126").Value = -32537.429
$ws.Range("H136").Value = 5560.2915
$ws.Range("I136").Value = 4726.2
$ws.Range("K136").Value = 14178.6
$ws.Range("M136").Value = -11628.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3686.9412
$ws.Range("I122").Value = 3147.9
$ws.Range("J122").Value = 4457
$ws.Range("K122").Value = 9443.700000000001
$ws.Range("L122").Value = 13371
$ws.Range("M122").Value = -6993.700000000001
$ws.Range("N122").Value = -18271
$ws.Range("H123").Value = 79666.336
$ws.Range("J123").Value = 79666.336
$ws.Range("L123").Value = 79666.336
$ws.Range("N123").Value = -89466.336
$ws.Range("H126").Value = 2977.65
$ws.Range("I126").Value = 2065.8125
$ws.Range("J126").Value = 6625
$ws.Range("K126").Value = 6197.4375
$ws.Range("L126").Value = 19875
$ws.Range("M126").Value = -3727.4375
$ws.Range("N126").Value = -24815
$ws.Range("H128").Value = 99997.25
$ws.Range("J128").Value = 99997.25
$ws.Range("L128").Value = 99997.25
$ws.Range("N128").Value = -109957.25
$ws.Range("J139").Value = 200000
$ws.Range("L139").Value = 200000
$ws.Range("N139").Value = -210280

